$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G8").Value = '/food/sort/price'
$ws.Range("G9").Value = '/food/search'
$ws.Range("G10").Value = '/food/type'
$ws.Range("G14").Value = '/store/food'
$ws.Range("G15").Value = '/store/food/{id}'
$ws.Range("F16").Value = '{StoreId:101
 Name:"StoreName"
 Street:"Store Street"
 MobileNo:"99xxxxxxxx"
 City:"StoreCity"
 State:"StoreState"
 Pincode:"Store Pincode"
 }'
$ws.Range("H16").Value = 'POST'
$ws.Range("F18").Value = '{
 StoreId:"101"
 Name:"StoreName"
 Street:"Store Street"
 MobileNo:"99xxxxxxxx"
 City:"StoreCity"
 State:"StoreState"
 Pincode:"Store Pincode"
 }'
$ws.Range("G18").Value = '/store'
$ws.Range("G22").Value = '/user/login'
$ws.Range("F23").Value = '{UserID:101
 "FirstName": "somename"
 "LastName": "somename"
 "DateOfBirth": "date"
 "Gender": "gender"
 "Street": "streetname"
 "Location": "landmark"
 "City": "cityname"
 "State": "statename"
 "Pincode": "pin" 
 "MobileNo": 9xxxxxxxx
 "EmailId": "abc@xyz.com"
"Password":"password"
}
'
$ws.Range("G23").Value = '/user
'
$ws.Range("H23").Value = 'POST'
$ws.Range("G24").Value = '/user/logout'
$ws.Range("G25").Value = '/user/profile'
$ws.Range("G26").Value = '/user/changepassword'
$ws.Range("F27").Value = '{
 UserId:"1001",
Fistname:"first name",
Lastname:"last name",
DateOfBirth:"0000-00-00",
Gender:"male/female",
Street:"street...",
Location:"location..",
City:"city",
User_State:"state",
Pincode:"110001",
MobileNo:"9985669826",
EmailId:"emailid@abc.com"
}'
$ws.Range("G27").Value = '/user'
$ws.Range("G28").Value = '/user/cart/item'
$ws.Range("H28").Value = 'POST'
$ws.Range("G29").Value = '/user/cart/item'
$ws.Range("F30").Value = '{OrderId:101
OrderDate:"YYYY-DD-MM"
StoreId:"storeId"
OrderStatus:"Confirmed"
cartid:"id of cart"
Street: "streetname"
Location: "landmark"
City: "cityname"
State: "statename"
PinCode:"pncode"
mobile number:"99xxxxxxxx"
totalprice:1500
}'
$ws.Range("G30").Value = '/food/order'
$ws.Range("G31").Value = '/user/order'
$ws.Range("G32").Value = '/user/order/filter/date'
$ws.Range("G33").Value = '/user/order/filter/status'
$ws.Range("F34").Value = '{
CreditCardNumber:"creditCardNumber"
 ValidFrom:"ValidFromDate in the form {month/year}"
 Valid to: "ValidToDate in the form {month/year}"
 Balance: Initial Balance in the account
 }
'
$ws.Range("G34").Value = '/user/creditcard'
$ws.Range("F35").Value = '{
CreditCardNumber:"creditCardNumber"
 ValidFrom:"ValidFromDate in the form {month/year}"
 Valid to: "ValidToDate in the form {month/year}"
 Balance: Initial Balance in the account
 }
'
$ws.Range("G35").Value = '/user/creditcard/verify'
$ws.Range("H35").Value = 'GET'
$ws.Range("G36").Value = '/user/creditcard'
$ws.Range("G37").Value = '/user/creditcard'
$ws.Range("G38").Value = '/user/cart'
$ws.Range("G40").Value = '/store/food'
$ws.Range("G40").WrapText = $true
$ws.Range("H40").Value = 'POST'
$ws.Range("H40").WrapText = $true
$ws.Range("G41").Value = '/store/food'
$ws.Range("G41").WrapText = $true
$ws.Range("H41").Value = 'DELETE'
$ws.Range("H41").WrapText = $true

# Row 35 height changed (content now wraps into more lines)
$ws.Rows(35).RowHeight = 165

# Update the active selection / scroll position to match the final view state
$ws.Range("H41").Select()
try { $excel.ActiveWindow.ScrollRow = 26 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
